$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D receive numeric-looking strings (e.g. "582.87").
# Force text entry so Excel does not auto-convert them to numbers,
# then restore the default "Normal" style so no residual number
# formatting is left behind on the cell.
$dCells = @("D2","D3","D5","D6","D7","D12","D14","D15","D17","D19","D20","D23","D26","D27","D30","D31","D32","D33","D36","D38","D39","D40","D42","D43","D44","D45","D47","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "68.133.96"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.256.46"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "582.87"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "185.06"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "3.831.14"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "28.21"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "68.233.60"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "3.244.15"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "393.20"
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "71.36"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "0.189"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").Value = "9.80"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "5.70"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "22.93"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "7.17"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").Value = "162.59"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +5.79%  "
$ws.Range("D38").Value = "0.823"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "26.88"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("D42").Value = "2.49"
$ws.Range("E42").Value = "  -5.28%  "
$ws.Range("D43").Value = "0.0691"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").Value = "2.654.58"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "25.40"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "339.59"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("D50").Value = "31.54"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").Value = "0.992"
$ws.Range("E51").Value = "  -1.01%  "

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }
